$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "BN-GC-14-1-o-felt-d"
$ws.Range("B2").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\BN-GC-14-1-o-felt-d\2_image_BN-GC-14-1-o-felt-d.jpg"
$ws.Range("C2").Value = 699
$ws.Range("D2").Value = 500
$ws.Range("E2").Value = "min side 500 < 501"

# Row 3
$ws.Range("A3").Value = "BN-OP-12-g-kr"
$ws.Range("B3").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\BN-OP-12-g-kr\001_image_BN-OP-12-g-kr.jpg"
$ws.Range("C3").Value = 620
$ws.Range("D3").Value = 500
$ws.Range("E3").Value = "min side 500 < 501"

# Row 4
$ws.Range("A4").Value = "BN-SB-6_073926493725"
$ws.Range("B4").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\BN-SB-6\3_image_BN-SB-6_073926493725.jpg"
$ws.Range("C4").Value = 643
$ws.Range("D4").Value = 500
$ws.Range("E4").Value = "min side 500 < 501"

# Row 5
$ws.Range("A5").Value = "TW-PH-black-ksr"
$ws.Range("B5").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\TW-PH-black-ksr\5_image_TW-PH-black-ksr.jpg"
$ws.Range("C5").Value = 500
$ws.Range("D5").Value = 716
$ws.Range("E5").Value = "min side 500 < 501"

# Row 6
$ws.Range("A6").Value = "TW-PH-kon-crz"
$ws.Range("B6").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\TW-PH-kon-crz\003_image_TW-PH-kon-crz.jpg"
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = 686
$ws.Range("E6").Value = "min side 500 < 501"

# Row 7
$ws.Range("A7").Value = "TW-PH-mars-ksr"
$ws.Range("B7").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\TW-PH-mars-ksr\003_image_TW-PH-mars-ksr.jpg"
$ws.Range("C7").Value = 500
$ws.Range("D7").Value = 716
$ws.Range("E7").Value = "min side 500 < 501"

# Row 8
$ws.Range("A8").Value = "TW-PH-red-saf"
$ws.Range("B8").Value = "C:/Users/Asus/Desktop/не всі фото/1 без фоторум\TW-PH-red-saf\001_image_TW-PH-red-saf.jpg"
$ws.Range("C8").Value = 500
$ws.Range("D8").Value = 645
$ws.Range("E8").Value = "min side 500 < 501"
